$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing "Sheet1" (the old beads-in-FSW data) to
#        "VvsHdiffVol_nogood" ------------------------------------------------
$oldWs = $wb.Worksheets.Item(1)
$oldWs.Name = "VvsHdiffVol_nogood"

# Select the whole first row on the old sheet (matches the recorded
# selection sqref="A1:XFD1" left behind on that tab) before we move away
# from it.
$oldWs.Rows.Item(1).Select() | Out-Null

# --- 2. Insert a brand-new sheet right after it; Excel names it "Sheet1"
#        automatically since that name just became free -------------------
$newWs = $wb.Worksheets.Add($null, $oldWs)

# --- 3. Populate the new "Sheet1" with the new vertical-vs-horizontal
#        9um-beads experiment data. Cells are written in the same order the
#        original author must have entered them (so newly-introduced shared
#        strings land at the same table positions) ---------------------------
$newWs.Range("E2").Value = "9um beads"
$newWs.Range("F2").Value = "2nd exp first file after high conc runs of last exp"
$newWs.Range("F1").Value = "Comments2"
$newWs.Range("A3").Value = "D20151103T150629"
$newWs.Range("A2").Value = "D20151103T145610"
$newWs.Range("A4").Value = "D20151103T151633"
$newWs.Range("A5").Value = "D20151103T152713"
$newWs.Range("F4").Value = "realize incorrectly running 2ml when want to run 5ml"
$newWs.Range("A6").Value = "D20151103T154933"
$newWs.Range("B4").Value = "2?"

# Header row (reuse existing shared strings from the old sheet)
$newWs.Range("A1").Value = "Filename"
$newWs.Range("B1").Value = "Volume"
$newWs.Range("C1").Value = "HorzOrVert"
$newWs.Range("D1").Value = "CellConc"
$newWs.Range("E1").Value = "Comments"

# Data rows
$newWs.Range("B2").Value = 2
$newWs.Range("C2").Value = "V"
$newWs.Range("D2").Value = 621

$newWs.Range("B3").Value = 2
$newWs.Range("C3").Value = "V"
$newWs.Range("D3").Value = 585
$newWs.Range("E3").Value = "9um beads"

$newWs.Range("C4").Value = "V"
$newWs.Range("D4").Value = 583
$newWs.Range("E4").Value = "9um beads"

$newWs.Range("B5").Value = 5
$newWs.Range("C5").Value = "V"
$newWs.Range("D5").Value = 633
$newWs.Range("E5").Value = "9um beads"

$newWs.Range("B6").Value = 5
$newWs.Range("C6").Value = "V"
$newWs.Range("E6").Value = "9um beads"
# (no D6 - left blank on purpose)

# --- 4. Copy the header cell's number format / alignment (centered, 0.00)
#        from the old sheet's D1 so the new D1 reuses the very same style --
$oldWs.Range("D1").Copy() | Out-Null
$newWs.Range("D1").PasteSpecial(-4122) | Out-Null

# --- 5. Column widths (best-fit-ish, matching the autofit widths Excel
#        would have computed for this content) ------------------------------
$newWs.Columns.Item(1).ColumnWidth = 16.736979166666668
$newWs.Columns.Item(2).ColumnWidth = 7.166666666666667
$newWs.Columns.Item(3).ColumnWidth = 10.166666666666666
$newWs.Columns.Item(4).ColumnWidth = 7.877604166666667
$newWs.Columns.Item(5).ColumnWidth = 9.736979166666666
$newWs.Columns.Item(6).ColumnWidth = 43.022135416666664

# --- 6. Final selection on the new sheet (matches recorded activeCell D6) --
$newWs.Range("D6").Select() | Out-Null
